$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-5 from 45185 (2023-09-16) to 45204 (2023-10-05)
$ws.Range("C2").Value = 45204
$ws.Range("C3").Value = 45204
$ws.Range("C4").Value = 45204
$ws.Range("C5").Value = 45204
